$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.554.39'
$ws.Range("E2").Value = '  -3.36%  '

# Row 3
$ws.Range("D3").Value = '1.849.11'
$ws.Range("E3").Value = '  -3.84%  '

# Row 4
$ws.Range("E4").Value = '  -0.89%  '

# Row 5
$ws.Range("D5").Value = '333.55'
$ws.Range("E5").Value = '  +2.36%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.84%  '

# Row 7
$ws.Range("D7").Value = '0.4661'
$ws.Range("E7").Value = '  -3.37%  '

# Row 8
$ws.Range("D8").Value = '0.3922'
$ws.Range("E8").Value = '  -3.68%  '

# Row 9
$ws.Range("D9").Value = '46.51'
$ws.Range("E9").Value = '  -2.44%  '

# Row 10
$ws.Range("D10").Value = '0.07903'
$ws.Range("E10").Value = '  -4.20%  '

# Row 11
$ws.Range("D11").Value = '0.9852'
$ws.Range("E11").Value = '  -2.66%  '

# Row 12
$ws.Range("E12").Value = '  -5.66%  '

# Row 13
$ws.Range("D13").Value = '1.973.80'
$ws.Range("E13").Value = '  +6.70%  '

# Row 14
$ws.Range("E14").Value = '  -3.81%  '

# Row 15
$ws.Range("D15").Value = '7.027'
$ws.Range("E15").Value = '  -3.42%  '

# Row 16
$ws.Range("D16").Value = '0.06889'
$ws.Range("E16").Value = '  +0.34%  '

# Row 17
$ws.Range("D17").Value = '87.72'
$ws.Range("E17").Value = '  -4.33%  '

# Row 18
$ws.Range("E18").Value = '  -0.82%  '

# Row 19
$ws.Range("D19").Value = '0.00001007'
$ws.Range("E19").Value = '  -3.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.10'
$ws.Range("E20").Value = '  -2.96%  '

# Row 21
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.71%  '

# Row 22
$ws.Range("D22").Value = '28.592.60'
$ws.Range("E22").Value = '  -3.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.400'
$ws.Range("E23").Value = '  -5.02%  '

# Row 24
$ws.Range("E24").Value = '  -5.17%  '

# Row 25
$ws.Range("D25").Value = '2.224.95'
$ws.Range("E25").Value = '  +6.32%  '

# Row 26
$ws.Range("D26").Value = '2.125'
$ws.Range("E26").Value = '  -2.71%  '

# Row 27
$ws.Range("D27").Value = '153.49'
$ws.Range("E27").Value = '  -1.63%  '

# Row 28
$ws.Range("D28").Value = '19.42'
$ws.Range("E28").Value = '  -2.89%  '

# Row 29
$ws.Range("D29").Value = '6.141'
$ws.Range("E29").Value = '  -5.53%  '

# Row 30
$ws.Range("D30").Value = '2.015'
$ws.Range("E30").Value = '  -4.11%  '

# Row 31
$ws.Range("D31").Value = '117.59'
$ws.Range("E31").Value = '  -2.61%  '

# Row 32
$ws.Range("D32").Value = '0.9845'
$ws.Range("E32").Value = '  -3.54%  '

# Row 33
$ws.Range("D33").Value = '0.09429'
$ws.Range("E33").Value = '  -2.47%  '

# Row 34
$ws.Range("D34").Value = '5.374'
$ws.Range("E34").Value = '  -4.69%  '

# Row 35
$ws.Range("D35").Value = '3.492'
$ws.Range("E35").Value = '  -1.72%  '

# Row 36
$ws.Range("E36").Value = '  -2.15%  '

# Row 37
$ws.Range("D37").Value = '0.06156'
$ws.Range("E37").Value = '  -3.69%  '

# Row 38
$ws.Range("D38").Value = '0.02205'
$ws.Range("E38").Value = '  -4.27%  '

# Row 39
$ws.Range("E39").Value = '  -2.08%  '

# Row 40
$ws.Range("D40").Value = '0.5713'

# Row 41
$ws.Range("D41").Value = '7.606'
$ws.Range("E41").Value = '  -3.86%  '

# Row 42
$ws.Range("D42").Value = '10.13'
$ws.Range("E42").Value = '  -5.95%  '

# Row 43
$ws.Range("D43").Value = '0.1795'
$ws.Range("E43").Value = '  -3.14%  '

# Row 44
$ws.Range("D44").Value = '2.368'
$ws.Range("E44").Value = '  -4.80%  '

# Row 45
$ws.Range("D45").Value = '1.252'
$ws.Range("E45").Value = '  +0.58%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5397'
$ws.Range("E46").Value = '  -3.24%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '11.82'
$ws.Range("E47").Value = '  -4.45%  '

# Row 48
$ws.Range("D48").Value = '0.07151'
$ws.Range("E48").Value = '  -4.73%  '

# Row 49
$ws.Range("D49").Value = '1.906'
$ws.Range("E49").Value = '  -2.42%  '

# Row 50
$ws.Range("D50").Value = '114.16'

# Row 51
$ws.Range("D51").Value = '42.98'
$ws.Range("E51").Value = '  +2.25%  '
